$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: paragraph under "Difficulty assessment" reads "No encountered
# difficulties" (no trailing period). Add a new run containing "." so the
# sentence ends consistently with its sibling paragraph (which already has
# the trailing period split into its own run).
# ---------------------------------------------------------------------------
$target1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq "No encountered difficulties") {
        $target1 = $p
    }
}

if ($target1 -ne $null) {
    $r = $target1.Range
    $xmlFragment = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4B869972" w14:textId="5A4A3DD4" w:rsidR="002A55D0" w:rsidRDefault="00607217"><w:pPr><w:rPr><w:rFonts w:ascii="Franklin Gothic Book" w:hAnsi="Franklin Gothic Book"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Franklin Gothic Book" w:hAnsi="Franklin Gothic Book"/></w:rPr><w:t>No encountered difficulties</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Franklin Gothic Book" w:hAnsi="Franklin Gothic Book"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
    [void]$r.InsertXML($xmlFragment)
}

# ---------------------------------------------------------------------------
# Edit 2: paragraph under "Improvement/Reflection" ends with two extra
# whitespace-only runs (" " and "  ") trailing the real sentence. Strip that
# trailing whitespace back off, leaving the sentence's own run untouched.
# ---------------------------------------------------------------------------
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "To improve this having the input being prompted again once the program has been run, would be a great addition.*") {
        $target2 = $p
    }
}

if ($target2 -ne $null) {
    $r2 = $target2.Range
    $sentence = "To improve this having the input being prompted again once the program has been run, would be a great addition."
    $delStart = $r2.Start + $sentence.Length
    $delEnd = $r2.End - 1
    if ($delEnd -gt $delStart) {
        $delRange = $d.Range($delStart, $delEnd)
        $delRange.Delete()
    }
}
